$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 11:22:58"
$wsZhCn.Range("G5").Value = "2016-01-25 11:23:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 11:23:08"
$wsDeDe.Range("G5").Value = "2016-01-25 11:23:54"
